$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 129, shifting existing rows 129:183 down to 130:184
$ws.Rows.Item(129).EntireRow.Insert()

# Populate the newly inserted row 129 with the new data point
$ws.Cells.Item(129, 1).Value = 3
$ws.Cells.Item(129, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(129, 3).Value = "Coquimbo"
$ws.Cells.Item(129, 4).Value = 44523
$ws.Cells.Item(129, 5).Value = 5
$ws.Cells.Item(129, 6).Value = 100112001
$ws.Cells.Item(129, 7).Value = "Berenjena"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 110
$ws.Cells.Item(129, 11).Value = 7500
$ws.Cells.Item(129, 12).Value = 8000
$ws.Cells.Item(129, 13).Value = 7727
$ws.Cells.Item(129, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(129, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(129, 16).Value = 129
$ws.Cells.Item(129, 17).Value = 60
$ws.Cells.Item(129, 18).Value = "Hortaliza"
